# Updated cryptos list on Sat Sep  9 22:42:25 UTC 2023 with GitHub Actions
# Refresh Price (col D) and Volume(1h) (col E) figures for the coin rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Some "Price" values (e.g. "215.38", "1.00") parse as plain numbers,
    # which would make Excel silently convert them from text to a numeric
    # cell. Force a text number format so the literal string is preserved,
    # matching the original "inline string" price column.
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.017.11"
$ws.Range("E2").Value = "  +0.22%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.643.28"
$ws.Range("E3").Value = "  +0.29%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.29%  "

# Row 5 - BNB
Set-TextValue "D5" "215.38"
$ws.Range("E5").Value = "  +0.28%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.03%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.28%  "

# Row 8 - Dogecoin
Set-TextValue "D8" "0.0638"
$ws.Range("E8").Value = "  +0.29%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.19%  "

# Row 10 - Solana
$ws.Range("E10").Value = "  -0.26%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.36%  "

# Row 12 - Polkadot
$ws.Range("E12").Value = "  +0.15%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.597.21"
$ws.Range("E13").Value = "  -1.90%  "

# Row 15 - Litecoin
$ws.Range("E15").Value = "  +1.40%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +0.44%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.058.02"
$ws.Range("E17").Value = "  +0.36%  "

# Row 18 - Dai
Set-TextValue "D18" "1.00"
$ws.Range("E18").Value = "  +0.27%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "194.50"
$ws.Range("E19").Value = "  +0.35%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  -0.41%  "

# Row 21 - Avalanche
$ws.Range("E21").Value = "  -0.33%  "

# Row 22 - Chainlink
$ws.Range("E22").Value = "  -1.01%  "

# Row 23 - Stellar
$ws.Range("E23").Value = "  +4.68%  "

# Row 24 - Monero
Set-TextValue "D24" "144.01"
$ws.Range("E24").Value = "  -0.17%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  -0.12%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.00%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  +0.50%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "15.51"
$ws.Range("E28").Value = "  +0.13%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  +0.43%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -1.02%  "

# Row 31 - Filecoin
$ws.Range("E31").Value = "  +0.90%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  -0.85%  "

# Row 33 - LidoDAOToken
$ws.Range("E33").Value = "  -0.21%  "

# Row 34 - HuobiToken
$ws.Range("E34").Value = "  +0.96%  "

# Row 35 - ARBITRUM
$ws.Range("E35").Value = "  +0.04%  "

# Row 36 - Maker
$ws.Range("D36").Value = "1.130.52"
$ws.Range("E36").Value = "  -0.82%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  -1.29%  "

# Row 38 - MXToken
$ws.Range("E38").Value = "  +0.06%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -0.18%  "

# Row 40 - FraxShare
$ws.Range("E40").Value = "  +0.51%  "

# Row 41 - Quant
Set-TextValue "D41" "98.79"
$ws.Range("E41").Value = "  -0.61%  "

# Row 42 - TrustWalletToken
Set-TextValue "D42" "0.795"
$ws.Range("E42").Value = "  -0.66%  "

# Row 43 - BabyDogeCoin
$ws.Range("E43").Value = "  +1.20%  "

# Row 44 - Aave
Set-TextValue "D44" "56.49"
$ws.Range("E44").Value = "  +0.01%  "

# Row 45 - RenderToken
Set-TextValue "D45" "1.50"
$ws.Range("E45").Value = "  +2.46%  "

# Row 46 - Cronos
Set-TextValue "D46" "0.0521"
$ws.Range("E46").Value = "  -1.50%  "

# Row 47 - EnergySwap
$ws.Range("E47").Value = "  +1.58%  "

# Row 48 - Mantle
$ws.Range("E48").Value = "  -0.28%  "

# Row 49 - USDD
$ws.Range("E49").Value = "  +0.16%  "

# Row 51 - Aptos
Set-TextValue "D51" "5.53"
$ws.Range("E51").Value = "  -0.04%  "
